$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Text label updates ---
$ws.Range("A12").Value = "Real estate (estimate)"
$ws.Range("A15").Value = "Bonds (actual through april 29)"
$ws.Range("A18").Value = "Stocks (actual through april 29)"

# --- Updated actual figures ---
$ws.Range("B6").Value = -4.1
$ws.Range("B9").Value = -2.2
$ws.Range("B12").Value = -2
$ws.Range("C12").Value = 1

# --- Replace formulas with plain estimated values ---
$ws.Range("B15").Value = 2.4
$ws.Range("C15").Value = 2.4
$ws.Range("B18").Value = 17
$ws.Range("C18").Value = 17

# --- Clear the SPX rows (28 & 29 source data / formulas) ---
$ws.Range("A28:C28").ClearContents()
$ws.Range("A29").ClearContents()
$ws.Range("B29:C29").ClearContents()

# --- Clear selection so reopened file doesn't pin A6 ---
$ws.Range("A1").Select()
